$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

$ws.Range("E12").Value = 0.3629913849272427
$ws.Range("E13").Value = 0.30866755954535602
$ws.Range("C39").Value = "e_CH18-220"
$ws.Range("C56").Value = "e_CH18-220"
$ws.Range("C67").Value = "e_w234983117-220"
$ws.Range("C122").Value = "e_w234983117-220"
$ws.Range("C125").Value = "e_w234983117-220"
$ws.Range("E150").Value = 0.21381383751804844
$ws.Range("E151").Value = 0.17206733071733712
$ws.Range("E152").Value = 0.18231505170803797
$ws.Range("E153").Value = 0.20988535532947597
$ws.Range("E154").Value = 0.1701049810444224
$ws.Range("E155").Value = 0.19532613932247714
$ws.Range("E156").Value = 0.16209575724687297
$ws.Range("E158").Value = 0.16085025627375071
$ws.Range("E161").Value = 0.13652468601509371
$ws.Range("E162").Value = 0.20006982412215921
$ws.Range("E163").Value = 0.15226887751132734
$ws.Range("E164").Value = 0.16485344960649678
$ws.Range("E165").Value = 0.15456128021356608
$ws.Range("E166").Value = 0.16629376698088194
$ws.Range("E167").Value = 0.15273795001145538
$ws.Range("E168").Value = 0.21640319337561012
$ws.Range("E169").Value = 0.19247860444770779
$ws.Range("E170").Value = 0.19228757088918788
$ws.Range("E171").Value = 0.19699531868281184
$ws.Range("E172").Value = 0.21063530390326943
$ws.Range("E173").Value = 0.19745398836539674
$ws.Range("E174").Value = 0.13549669849969209
$ws.Range("C175").Value = "elc_spv-CHE_0010"
$ws.Range("C176").Value = "elc_spv-CHE_0003"
$ws.Range("C177").Value = "elc_spv-CHE_0003"
$ws.Range("C178").Value = "elc_spv-CHE_0003"
$ws.Range("C179").Value = "elc_spv-CHE_0003"
$ws.Range("C180").Value = "elc_spv-CHE_0021"
$ws.Range("C181").Value = "elc_spv-CHE_0021"
$ws.Range("C182").Value = "elc_spv-CHE_0007"
$ws.Range("C183").Value = "elc_spv-CHE_0007"
$ws.Range("C184").Value = "elc_spv-CHE_0007"
$ws.Range("C185").Value = "elc_spv-CHE_0007"
$ws.Range("C186").Value = "elc_spv-CHE_0006"
$ws.Range("C187").Value = "elc_spv-CHE_0006"
$ws.Range("C188").Value = "elc_spv-CHE_0006"
$ws.Range("C189").Value = "elc_spv-CHE_0000"
$ws.Range("C190").Value = "elc_spv-CHE_0000"
$ws.Range("C191").Value = "elc_spv-CHE_0019"
$ws.Range("C192").Value = "elc_spv-CHE_0019"
$ws.Range("C193").Value = "elc_spv-CHE_0019"
$ws.Range("C194").Value = "elc_spv-CHE_0018"
$ws.Range("C195").Value = "elc_spv-CHE_0018"
$ws.Range("C196").Value = "elc_spv-CHE_0018"
$ws.Range("C197").Value = "elc_spv-CHE_0018"
$ws.Range("C198").Value = "elc_spv-CHE_0018"
$ws.Range("C199").Value = "elc_spv-CHE_0025"
$ws.Range("C200").Value = "elc_spv-CHE_0025"
$ws.Range("C201").Value = "elc_spv-CHE_0025"
$ws.Range("C202").Value = "elc_spv-CHE_0025"
$ws.Range("C203").Value = "elc_spv-CHE_0011"
$ws.Range("C204").Value = "elc_spv-CHE_0011"
$ws.Range("C205").Value = "elc_spv-CHE_0011"
$ws.Range("C206").Value = "elc_spv-CHE_0011"
$ws.Range("C207").Value = "elc_spv-CHE_0011"
$ws.Range("C208").Value = "elc_spv-CHE_0017"
$ws.Range("C209").Value = "elc_spv-CHE_0022"
$ws.Range("C210").Value = "elc_spv-CHE_0002"
$ws.Range("C211").Value = "elc_spv-CHE_0002"
$ws.Range("C212").Value = "elc_spv-CHE_0002"
$ws.Range("C213").Value = "elc_spv-CHE_0002"
$ws.Range("C214").Value = "elc_spv-CHE_0009"
$ws.Range("C215").Value = "elc_spv-CHE_0009"
$ws.Range("C216").Value = "elc_spv-CHE_0015"
